$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")
$nl = [char]10

# 1. Rename the "Config" sheet to "Configv1" (hdpv1 configs)
$ws.Name = "Configv1"

# 2. Update the cell texts that now reference "dfs.replication" instead of the
#    generic "Replication Factor" wording, and extend the JT URI cell with the
#    two extra mapred task properties. The order below matches the order the
#    shared-strings table ends up in after the edit.
$ws.Range("G8").Value = "JT URI" + $nl + "mapred.hosts.exclude" + $nl + "mapred.hosts" + $nl + "mapred.map.tasks.maximum" + $nl + "mapred.reduce.tasks.maximum"
$ws.Range("K8").Value = "JT URI"
$ws.Range("C6").Value = "NN Dir" + $nl + "Block size" + $nl + "dfs.replication" + $nl + "dfs.hosts.exclude" + $nl + "dfs.hosts"
$ws.Range("G6").Value = "NA" + $nl + "dfs.replication"
$ws.Range("I6").Value = "DN Dir" + $nl + "dfs.replication"

# 3. Row 8 grew from 3 lines to 5 lines of wrapped text, so its height doubled.
$ws.Rows.Item(8).RowHeight = 105

# 4. Columns E and G were narrowed/widened to fit the new (shorter / longer)
#    text after the rewrite.
$ws.Columns.Item(5).ColumnWidth = 28.166666666666668
$ws.Columns.Item(7).ColumnWidth = 28.833333333333332

# 5. Update the remembered selection to the new last-edited cell.
$ws.Range("I6").Select() | Out-Null
